$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("Katian" stage) was originally using the abbreviated shared string
# "Kat." in column C (same as column B); correct it to the full stage name
# "Katian" (matching column A) as part of updating the 2014 conodont Sr data.
$ws.Range("C3").Value = "Katian"

# Move the current selection/active cell to C3, matching the saved view state.
$ws.Range("C3").Select()

# Reflect the updated workbook window position from the saved view state.
$aw = $excel.ActiveWindow
$aw.WindowState = -4143
$aw.Left = 3450
$aw.Top = 1830
